$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F9").Value = "En revisión editor"
$ws.Range("C1").Value = "Ciencias Naturales"
$ws.Range("C2").Value = 5
$ws.Range("C2:D2").Select()
